$d = $word.ActiveDocument

# Locate the paragraph "A Cat, A Parrot, and a Bag of Seed" so the new
# content lands right after the problem's title, before the blank line
# that precedes "Socks in the Dark".
$titleIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    $t = $t.Replace("`r", "").Replace("`a", "")
    if ($t -eq "A Cat, A Parrot, and a Bag of Seed") {
        $titleIndex = $i
        break
    }
}
if ($titleIndex -eq 0) {
    throw "Could not find paragraph 'A Cat, A Parrot, and a Bag of Seed'"
}

# Insert a blank paragraph right after the title.
$titlePara = $d.Paragraphs.Item($titleIndex)
$titlePara.Range.InsertParagraphAfter()

# Insert the "1) Define the Problem" heading paragraph after the new blank one.
$blankPara = $d.Paragraphs.Item($titleIndex + 1)
$blankPara.Range.InsertParagraphAfter()
$headingPara = $d.Paragraphs.Item($titleIndex + 2)
$headingPara.Range.Text = "1) Define the Problem"

# Insert the explanatory paragraph after the heading paragraph.
$headingPara.Range.InsertParagraphAfter()
$bodyPara = $d.Paragraphs.Item($titleIndex + 3)
$bodyPara.Range.Text = "In this scenario, the man must transport three things across a river in a boat that can only carry himself and one of the things. Meaning that he has to take them one at a time. The added catch is that, when left alone together, some of the things will destroy other things. So the man has to figure out what order to take the things safely across in."
